$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 19: Krashtkid "Spit" special attack entry
$ws.Range("A19").Value = "Spit"
$ws.Range("B19").Value = 80
$ws.Range("F19").Value = 1
$ws.Range("H19").Value = "Krashtkid"
$ws.Range("I19").Value = "Spit"
$ws.Range("J19").Value = "Spit: The spit of the krarshtkid is a gummy, web-like mess called pratzim. It has a potency of 15. A victim hit is entangled by the saliva and must make a STR resistance roll against the potency of the pratzim. A missed roll means that the victim is bound, may not fight with any weapon, and may move only at a rate of 2. Each subsequent round the victim may attempt another resistance roll to overcome the spit’s effect. It will remain on the adventurer until all garments can be thoroughly cleansed"

# Column widths widened (H to fit "Krashtkid" header-ish content, J to fit the long text)
$ws.Columns.Item(8).ColumnWidth = 9.857142857142858
$ws.Columns.Item(10).ColumnWidth = 457.57142857142856
